$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the airport/vertiport labels in column A
$ws.Range("A2").Value = "Vertiport 1"
$ws.Range("A6").Value = "Vertiport 2"
$ws.Range("A10").Value = "Vertiport 3"

# Rename the pad labels in column C for vertiport 1 (rows 2-5)
$ws.Range("C2").Value = "Vertiport1_Pad1"
$ws.Range("C3").Value = "Vertiport1_Pad2"
$ws.Range("C4").Value = "Vertiport1_Pad3"
$ws.Range("C5").Value = "Vertiport1_Pad4"

# Rename the pad labels in column C for vertiport 2 (rows 6-9)
$ws.Range("C6").Value = "Vertiport2_Pad1"
$ws.Range("C7").Value = "Vertiport2_Pad2"
$ws.Range("C8").Value = "Vertiport2_Pad3"
$ws.Range("C9").Value = "Vertiport2_Pad4"

# Rename the pad labels in column C for vertiport 3 (rows 10-13)
$ws.Range("C10").Value = "Vertiport3_Pad1"
$ws.Range("C11").Value = "Vertiport3_Pad2"
$ws.Range("C12").Value = "Vertiport3_Pad3"
$ws.Range("C13").Value = "Vertiport3_Pad4"

# Update the selection to match the saved view state
$ws.Range("A2:A10").Select()
